# ELI-114: fixing fixture script data that could lead to occasional test
# failures because it is invalid with respect to validation rules.
#
# The "(0)" trunk-prefix parenthetical notation is removed from the UK/NL
# phone numbers (e.g. "+44 (0)20 111 2222" -> "+44 20 111 2222") since that
# notation isn't a valid E.164-ish phone number for the importer's
# validation rules. A leading apostrophe is used (just like a real user
# typing in Excel) so the value is kept as literal text with the
# quote-prefix cell style preserved, rather than Excel trying to reinterpret
# the leading "+" as the start of a formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "'+44 20 111 2222"
$ws.Range("F2").Value = "'+44 7770 222 123"
$ws.Range("G2").Value = "'+44 1233 444 555"

$ws.Range("E3").Value = "'+31 20 999 1111"
$ws.Range("F3").Value = "'+31 6 4444 3333"

$ws.Range("E4").Value = "'+44 20 7777 3333"
$ws.Range("F4").Value = "'+44 7770 222 111"
$ws.Range("G4").Value = "'+44 1892 999 222"

$ws.Range("E5").Value = "'+31 20-555 1000"

$ws.Range("E6").Value = "'+31 20 222 1234"

$ws.Range("E7").Value = "'+31 20 456 7891"
$ws.Range("F7").Value = "'+31 6 432 1234"

$ws.Range("E8").Value = "'+31 20 444 1234"
$ws.Range("F8").Value = "'+31 6 112 2334"

# Reflect the cursor position left behind in the saved file.
$ws.Range("D15").Select()
